$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: update Correspond Handoff Datetime (D) and Correspond Handback DateTime (G)
# for rows 33 and 34 (same values on both rows, mirrored in B3822ce5... source file group)
$wsZhCn.Range("D33").Value = "2016-03-09 06:21:53"
$wsZhCn.Range("D34").Value = "2016-03-09 06:21:53"
$wsZhCn.Range("G33").Value = "2016-03-09 06:22:56"
$wsZhCn.Range("G34").Value = "2016-03-09 06:22:56"

# de-de sheet: same rows/columns
$wsDeDe.Range("D33").Value = "2016-03-09 06:21:57"
$wsDeDe.Range("D34").Value = "2016-03-09 06:21:57"
$wsDeDe.Range("G33").Value = "2016-03-09 06:23:03"
$wsDeDe.Range("G34").Value = "2016-03-09 06:23:03"
